$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A1").Value = "HK_R_acc_G"

$ws.Cells.Item(2, 1).Value = 93.023255813953483
$ws.Cells.Item(3, 1).Value = 92.389006342494724
$ws.Cells.Item(4, 1).Value = 93.128964059196619
$ws.Cells.Item(5, 1).Value = 93.446088794926013
$ws.Cells.Item(6, 1).Value = 93.446088794926013
$ws.Cells.Item(7, 1).Value = 96.40591966173362
$ws.Cells.Item(8, 1).Value = 93.128964059196619
$ws.Cells.Item(9, 1).Value = 93.128964059196619
$ws.Cells.Item(10, 1).Value = 92.494714587737846
$ws.Cells.Item(11, 1).Value = 92.494714587737846
$ws.Cells.Item(12, 1).Value = 95.560253699788589
$ws.Cells.Item(13, 1).Value = 95.454545454545453
$ws.Cells.Item(14, 1).Value = 93.128964059196619
$ws.Cells.Item(15, 1).Value = 93.128964059196619
$ws.Cells.Item(16, 1).Value = 93.023255813953483
$ws.Cells.Item(17, 1).Value = 93.446088794926013
$ws.Cells.Item(18, 1).Value = 96.40591966173362
$ws.Cells.Item(19, 1).Value = 96.40591966173362
$ws.Cells.Item(20, 1).Value = 93.657505285412256
$ws.Cells.Item(21, 1).Value = 93.657505285412256
$ws.Cells.Item(22, 1).Value = 93.763213530655392
$ws.Cells.Item(23, 1).Value = 95.877378435517969
$ws.Cells.Item(24, 1).Value = 95.877378435517969
$ws.Cells.Item(25, 1).Value = 95.560253699788589
$ws.Cells.Item(26, 1).Value = 92.811839323467225
$ws.Cells.Item(27, 1).Value = 92.811839323467225
$ws.Cells.Item(28, 1).Value = 92.811839323467225
$ws.Cells.Item(29, 1).Value = 93.023255813953483
$ws.Cells.Item(30, 1).Value = 93.023255813953483
$ws.Cells.Item(31, 1).Value = 93.023255813953483
$ws.Cells.Item(32, 1).Value = 92.917547568710361
$ws.Cells.Item(33, 1).Value = 93.551797040169134
$ws.Cells.Item(34, 1).Value = 92.917547568710361
$ws.Cells.Item(35, 1).Value = 93.128964059196619
$ws.Cells.Item(36, 1).Value = 93.234672304439741
$ws.Cells.Item(37, 1).Value = 96.511627906976756
$ws.Cells.Item(38, 1).Value = 93.234672304439741
$ws.Cells.Item(39, 1).Value = 93.023255813953483
$ws.Cells.Item(40, 1).Value = 95.983086680761105
$ws.Cells.Item(41, 1).Value = 92.811839323467225
$ws.Cells.Item(42, 1).Value = 92.811839323467225
$ws.Cells.Item(43, 1).Value = 92.811839323467225
$ws.Cells.Item(44, 1).Value = 92.811839323467225
$ws.Cells.Item(45, 1).Value = 92.706131078224104
$ws.Cells.Item(46, 1).Value = 92.811839323467225
$ws.Cells.Item(47, 1).Value = 92.177589852008452
$ws.Cells.Item(48, 1).Value = 92.811839323467225
$ws.Cells.Item(49, 1).Value = 92.177589852008452
